# Weekly refresh of the Fruta/Hortalizas "Ciruela" sheet.
# Each data row (2-22) is re-populated with the values that, in the
# source workbook, belonged to another row (rows were re-paired to new
# report dates while row 8 stays put). We snapshot the "before" values
# of the relevant columns for every row first, then write the
# "after" values, so the row-to-row copy isn't clobbered mid-flight.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a group for each record.
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S")

# target row -> source row (source row's CURRENT/"before" values become
# the target row's new values)
$mapping = @{
    2  = 17
    3  = 18
    4  = 11
    5  = 21
    6  = 5
    7  = 12
    8  = 8
    9  = 13
    10 = 16
    11 = 14
    12 = 7
    13 = 3
    14 = 2
    15 = 22
    16 = 6
    17 = 19
    18 = 15
    19 = 4
    20 = 9
    21 = 10
    22 = 20
}

# Snapshot current values for every row/column we might need as a source.
$snapshot = @{}
foreach ($row in 2..22) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowVals
}

# Apply the new values taken from the snapshot of the mapped source row.
foreach ($row in 2..22) {
    $srcRow = $mapping[$row]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $srcVals[$col]
    }
}
